# Burndown.xlsx update — sprint 5/6 data entry.
#
# Sprint 5 (row 7) now has a "completed" figure (46) and its computed
# "end of sprint" remainder (84 - 46 = 38).
# Sprint 6 (row 8) begins where sprint 5 left off (46).
# Sprint 7 (row 9) no longer has a stray placeholder "beginning" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = 46
$ws.Range("D7").Value = 38
$ws.Range("B8").Value = 46
$ws.Range("B9").ClearContents()

# The chart was nudged slightly (left/down) on the sheet — keep its size,
# just re-anchor its position to match.
$co = $ws.ChartObjects(1)
$co.Left = 260
$co.Top = 23.625
